$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block: account holder name ---
$ws.Range("C2").Value = "Hartmut"

# Card number (B3) is a 16-digit string; Excel's automatic type detection
# would otherwise coerce a plain digit string to a number, so force Text
# formatting first (mirrors how a real statement's card-number column is
# kept as text), then restore the cell's original formatting (copied from
# D3, which already carries the same style used by B3) so the cell keeps
# its original look-and-feel / style index.
$ws.Range("D3").Copy()
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("D3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 01.10.2023"

# --- Row 6 (existing transaction, values updated) ---
$ws.Range("B6").Value = "03.10."
$ws.Range("C6").Value = "04.10."
$ws.Range("D6").Value = "KARTENZ./03.10 ALDI SUED RO"
$ws.Range("E6").Value = "72,99-"

# --- Row 7 (existing transaction, values updated) ---
$ws.Range("B7").Value = "06.10."
$ws.Range("C7").Value = "07.10."
$ws.Range("D7").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E7").Value = "43,74-"

# --- Row 8 (existing transaction, values updated) ---
$ws.Range("B8").Value = "07.10."
$ws.Range("C8").Value = "08.10."
$ws.Range("D8").Value = "PAYPAL FUJGLE"
$ws.Range("E8").Value = "28,88-"

# --- Row 9 (previously blank, now a new transaction) ---
$ws.Range("B9").Value = "09.10."
$ws.Range("C9").Value = "10.10."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "24,97-"

# --- Row 10 (previously blank, now a new transaction) ---
$ws.Range("B10").Value = "12.10."
$ws.Range("C10").Value = "13.10."
$ws.Range("D10").Value = "PAYPAL DMEHUU"
$ws.Range("E10").Value = "60,73-"

# --- Row 11 (previously blank, now a new transaction) ---
$ws.Range("B11").Value = "16.10."
$ws.Range("C11").Value = "17.10."
$ws.Range("D11").Value = "KARTENZ./16.10 ALDI SUED RO"
$ws.Range("E11").Value = "141,34-"

# Amount column (E) for the newly-populated rows 9-11 used a different
# (wrapped / centered) number format before; bring it in line with the
# format already used by rows 6-8 (right aligned, no wrap) by copying it.
$ws.Range("E6").Copy()
$ws.Range("E9:E11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 18.10.2023"
$ws.Range("E12").Value = "372,65-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 24.10.2023"
